$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New relayConfig PGN row (row 18): label, byte count, and 20 numbered data
# bytes followed by the CRC column, matching the other PGN rows' layout.
$ws.Range("A18").Value = "relayConfig"
$ws.Range("E18").Value = 20

$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = 6
$ws.Range("L18").Value = 7
$ws.Range("M18").Value = 8
$ws.Range("N18").Value = 9
$ws.Range("O18").Value = 10
$ws.Range("P18").Value = 11
$ws.Range("Q18").Value = 12
$ws.Range("R18").Value = 13
$ws.Range("S18").Value = 14
$ws.Range("T18").Value = 15
$ws.Range("U18").Value = 16
$ws.Range("V18").Value = 17
$ws.Range("W18").Value = 18
$ws.Range("X18").Value = 19
$ws.Range("Y18").Value = 20
$ws.Range("Z18").Value = "CRC"

# Match the author's final selection / scroll position.
[void]$ws.Range("F19").Select()
